# Insert a new weekly price-report row at row 876 (shifting all subsequent
# rows down by one, dimension A1:R941 -> A1:R942) and populate it with the
# new "Poroto verde" / Magnum / Primera record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(876).EntireRow.Insert()

$ws.Cells.Item(876, 1).Value  = 6
$ws.Cells.Item(876, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(876, 3).Value  = "Metropolitana"
$ws.Cells.Item(876, 4).Value  = 44783
$ws.Cells.Item(876, 5).Value  = 13
$ws.Cells.Item(876, 6).Value  = 100112031
$ws.Cells.Item(876, 7).Value  = "Poroto verde"
$ws.Cells.Item(876, 8).Value  = "Magnum"
$ws.Cells.Item(876, 9).Value  = "Primera"
$ws.Cells.Item(876, 10).Value = 1000
$ws.Cells.Item(876, 11).Value = 30000
$ws.Cells.Item(876, 12).Value = 32000
$ws.Cells.Item(876, 13).Value = 31060
$ws.Cells.Item(876, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(876, 15).Value = "Perú"
$ws.Cells.Item(876, 16).Value = 1242
$ws.Cells.Item(876, 17).Value = 25
$ws.Cells.Item(876, 18).Value = "Hortaliza"
